# Add two new Test Case rows (TC_008, TC_009) to the "Test Cases" worksheet,
# mirroring the existing rows 10-16, and update the sheet selection/scroll
# position to reflect where the author ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$steps = '1. Launch Browser' + [char]10 + '2. Go to URL https://awesomeqa.com/ui/' + [char]10 + '3. Click on My Account' + [char]10 + '4. Click on Login'

# --- Row 17 (TC_008) : fill everything except the "Test Cases" (F) column first ---
$ws.Range("A17").Value = 'TC_008'
$ws.Range("B17").Value = 'Login '
$ws.Range("C17").Value = 'awsomeqa Login Page'
$ws.Range("D17").Value = $steps
$ws.Range("E17").Value = 'Email ID:- $#@&#@$'
$ws.Range("G17").Value = 'P0'
$ws.Range("H17").Value = 'Invalid Credentials'

# --- Row 18 (TC_009) : fill completely ---
$ws.Range("A18").Value = 'TC_009'
$ws.Range("B18").Value = 'Login '
$ws.Range("C18").Value = 'awsomeqa Login Page'
$ws.Range("D18").Value = $steps
$ws.Range("E18").Value = 'Password:- #@$%^%$#'
$ws.Range("F18").Value = 'Verify login with special characters ' + [char]10 + 'in password field'
$ws.Range("G18").Value = 'P0'
$ws.Range("H18").Value = 'Invalid Credentials'

# --- Row 17 column F filled last (matches original authoring/shared-string order) ---
$ws.Range("F17").Value = 'Verify login with special characters' + [char]10 + 'in Email ID field'

# Match formatting used by the other data rows: wrap text on columns D:H,
# and the custom row heights used for these two new rows.
$ws.Range("D17:H18").WrapText = $true
$ws.Rows.Item(17).RowHeight = 96
$ws.Rows.Item(18).RowHeight = 95.25

# Update the view: scroll down a bit and move the selection, as happened in
# the authored workbook (topLeftCell moved from A5 to A16, selection moved
# from H16 to G21).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G21").Select()
